$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new "Changed panorama" entry as row 5 first so it gets the
# earlier shared-string slot
$ws.Range("A5").Value = "Changed panorama"

# Update existing changelog entry text (row 4 -> "Fixed enchantment glint not working (1.21.2)")
$ws.Range("A4").Value = "Fixed enchantment glint not working (1.21.2)"

# Insert the new "Fixed VRAM leakage" entry as row 6
$ws.Range("A6").Value = "Fixed VRAM leakage"

# Copy the style from A4 to the new rows so formatting matches
$ws.Range("A4").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)

# Update selection to match the target state
$ws.Range("A7").Select()
